# Updates to status for planning next few weeks of development
#
# 1. Re-point Sheet1 / Sheet2 selections (no tab stays "selected" on them).
# 2. Add a new "Sheet3" ("MS Sched") with a day-by-day schedule table.
# 3. Make Sheet3 the active sheet with its own selection.

$wb = $excel.ActiveWorkbook

# --- Sheet1: drop the old selection, move it to D16 -----------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("D16").Select() | Out-Null

# --- Sheet2: move the selection to A10 -------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A10").Select() | Out-Null

# --- Add Sheet3 at the end of the tab strip --------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Title / "as at" header band
$ws3.Range("A1").Value = "MS Sched"
$ws3.Range("I1").Value = "As at 08/02/19"
$ws3.Range("I2").Value = "Waiting on"

$ws3.Range("I3").Value = "Rabbi Bordon"
$ws3.Range("J3").Value = "MD"

$ws3.Range("E4").Value = "Complete Jsilv"
$ws3.Range("I4").Value = "Rabbi Klein"
$ws3.Range("J4").Value = "Him"

$ws3.Range("E5").Value = "Type up Jsilv"
$ws3.Range("I5").Value = "MS"
$ws3.Range("J5").Value = "ML"

$ws3.Range("E6").Value  = "Page 1 MD"
$ws3.Range("E7").Value  = "Page 2 MD"
$ws3.Range("E8").Value  = "Page 3 MD"
$ws3.Range("E10").Value = "Page 4 MD"
$ws3.Range("E11").Value = "Complete MD"
$ws3.Range("E12").Value = "Computer MD"
$ws3.Range("E13").Value = "Section 1,2 SE"
$ws3.Range("E14").Value = "Section 3 SE"
$ws3.Range("E15").Value = "Section 4,5 SE"
$ws3.Range("E17").Value = "1/2 SE Computer"
$ws3.Range("E18").Value = "Final check of MK"
$ws3.Range("E19").Value = "2/2 SE Computer"
$ws3.Range("E20").Value = "1/2 Final check of DS"
$ws3.Range("E21").Value = "2/2 Final check of DS"

# Date column (C) + weekday formula column (D), rows 3-26
for ($r = 3; $r -le 26; $r++) {
    $c = $ws3.Range("C$r")
    $c.Value = 43504 + ($r - 3)
    $c.NumberFormat = "mm-dd-yy"
}

$ws3.Range("D3").Formula = '=IF(TEXT(C3, "dddd")="Saturday", "Motzash", TEXT(C3, "dddd"))'
$ws3.Range("D4:D26").FormulaR1C1 = '=IF(TEXT(RC[-1], "dddd")="Saturday", "Motzash", TEXT(RC[-1], "dddd"))'

# Column C is sized to fit the dates
$ws3.Columns.Item(3).ColumnWidth = 9.8

# Sheet3 becomes the active tab, selection resting on J6
$ws3.Range("J6").Select() | Out-Null
